$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - column F holds "想去人数" (interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1051
$ws1.Range("F3").Value = 318
$ws1.Range("F4").Value = 2848
$ws1.Range("F6").Value = 603

# Sheet "全部类型" (All types) - same values, offset by two rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1051
$ws4.Range("F5").Value = 318
$ws4.Range("F6").Value = 2848
$ws4.Range("F8").Value = 603
